# Update the "dSF" column (F) values for a subset of rows to reflect
# repulled/pushed data and corrected mean calculation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -1
    3  = -1
    13 = -2
    15 = 0
    21 = 4
    23 = -3
    26 = 0
    28 = -2
    38 = 4
    39 = -1
    40 = -1
    53 = 3
    54 = -2
    56 = -4
    57 = -3
    64 = 4
    67 = 0
    74 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
